$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A3").Value = "Pipeline(steps=[('scaler', RobustScaler()),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f91044dc400>),
                ('model',
                 BaggingClassifier(estimator=LogisticRegression(C=1,
                                                                class_weight='balanced',
                                                                max_iter=1000,
                                                                penalty='l1',
                                                                random_state=42,
                                                                solver='saga'),
                                   n_estimators=50, random_state=42))])"

$ws.Range("C3").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f91044d0dc0>, 'scaler': RobustScaler(), 'model__n_estimators': 50, 'model__estimator__solver': 'saga', 'model__estimator__penalty': 'l1', 'model__estimator__class_weight': 'balanced', 'model__estimator__C': 1}"

$ws.Range("A5").Value = "Pipeline(steps=[('scaler', StandardScaler()),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f91044dc490>),
                ('model',
                 BaggingClassifier(estimator=LogisticRegression(C=0.0001,
                                                                max_iter=1000,
                                                                random_state=42,
                                                                solver='liblinear'),
                                   random_state=42))])"

$ws.Range("C5").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f91c789b460>, 'scaler': StandardScaler(), 'model__n_estimators': 10, 'model__estimator__solver': 'liblinear', 'model__estimator__penalty': 'l2', 'model__estimator__class_weight': None, 'model__estimator__C': 0.0001}"
